$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12, shifting the existing rows 12-21 down to 13-22.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly entry.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44447
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = 100112044
$ws.Cells.Item(12, 7).Value = "Perejil"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 300
$ws.Cells.Item(12, 11).Value = 900
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 13).Value = 950
$ws.Cells.Item(12, 14).Value = "`$/atado 1,5 a 2 kilos"
$ws.Cells.Item(12, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(12, 16).Value = 475
$ws.Cells.Item(12, 17).Value = 2
$ws.Cells.Item(12, 18).Value = "Hortaliza"
